$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1062.86
$ws.Range("I15").Value = 1062.86
$ws.Range("K15").Value = 3188.58
$ws.Range("M15").Value = -3019.58
$ws.Range("H33").Value = 444.5
$ws.Range("J33").Value = 280.66666
$ws.Range("L33").Value = 280.66666
$ws.Range("N33").Value = -738.66666
$ws.Range("H74").Value = 5499.5
$ws.Range("I74").Value = 5499.5
$ws.Range("K74").Value = 5499.5
$ws.Range("M74").Value = -4563.5
$ws.Range("H77").Value = 5499.5
$ws.Range("I77").Value = 5499.5
$ws.Range("K77").Value = 27497.5
$ws.Range("M77").Value = -22817.5
$ws.Range("H96").Value = 2504.4285
$ws.Range("I96").Value = 2750.3333
$ws.Range("K96").Value = 8250.999899999999
$ws.Range("M96").Value = -6877.999899999999
$ws.Range("H129").Value = 879.38464
$ws.Range("J129").Value = 935.1111
$ws.Range("L129").Value = 2805.3333
$ws.Range("N129").Value = -12805.3333
$ws.Range("H137").Value = 4178.2856
$ws.Range("I137").Value = 4476.4
$ws.Range("K137").Value = 13429.2
$ws.Range("M137").Value = -10879.2
$ws.Range("H138").Value = 2897.3374
$ws.Range("I138").Value = 3150
$ws.Range("J138").Value = 2876.8513
$ws.Range("K138").Value = 9450
$ws.Range("L138").Value = 8630.553899999999
$ws.Range("M138").Value = -4310
$ws.Range("N138").Value = -18910.5539

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11198.3
$ws.Range("I32").Value = 8323.777
$ws.Range("K32").Value = 8323.777
$ws.Range("M32").Value = -8036.777
$ws.Range("H74").Value = 1690.5358
$ws.Range("I74").Value = 942.44446
$ws.Range("J74").Value = 3037.1
$ws.Range("K74").Value = 942.44446
$ws.Range("L74").Value = 3037.1
$ws.Range("M74").Value = -68.44446000000005
$ws.Range("N74").Value = -4785.1
$ws.Range("H77").Value = 1690.5358
$ws.Range("I77").Value = 942.44446
$ws.Range("J77").Value = 3037.1
$ws.Range("K77").Value = 4712.2223
$ws.Range("L77").Value = 15185.5
$ws.Range("M77").Value = -344.2223000000004
$ws.Range("N77").Value = -23921.5
$ws.Range("H110").Value = 569.1429000000001
$ws.Range("I110").Value = 536
$ws.Range("K110").Value = 536
$ws.Range("M110").Value = 1509
$ws.Range("H132").Value = 2732.8708
$ws.Range("I132").Value = 2610.0454
$ws.Range("J132").Value = 3033.111
$ws.Range("K132").Value = 7830.1362
$ws.Range("L132").Value = 9099.332999999999
$ws.Range("M132").Value = -5300.1362
$ws.Range("N132").Value = -14159.333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 7576458
$ws.Range("I94").Value = 9259925
$ws.Range("J94").Value = 856.6667
$ws.Range("K94").Value = 9259925
$ws.Range("L94").Value = 856.6667
$ws.Range("M94").Value = -9259474
$ws.Range("N94").Value = -1758.6667
$ws.Range("H107").Value = 940.8421
$ws.Range("I107").Value = 879.75
$ws.Range("J107").Value = 1266.6666
$ws.Range("K107").Value = 879.75
$ws.Range("L107").Value = 1266.6666
$ws.Range("M107").Value = 1040.25
$ws.Range("N107").Value = -5106.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2042.05
$ws.Range("I31").Value = 1966.4615
$ws.Range("K31").Value = 1966.4615
$ws.Range("M31").Value = -1671.4615
$ws.Range("H34").Value = 2042.05
$ws.Range("I34").Value = 1966.4615
$ws.Range("K34").Value = 1966.4615
$ws.Range("M34").Value = -1764.4615
$ws.Range("H107").Value = 893.3333
$ws.Range("I107").Value = 635
$ws.Range("J107").Value = 1410
$ws.Range("K107").Value = 635
$ws.Range("L107").Value = 1410
$ws.Range("M107").Value = 1285
$ws.Range("N107").Value = -5250
$ws.Range("H132").Value = 2195
$ws.Range("I132").Value = 1888.8
$ws.Range("J132").Value = 3045.5557
$ws.Range("K132").Value = 5666.4
$ws.Range("L132").Value = 9136.667099999999
$ws.Range("M132").Value = -3136.4
$ws.Range("N132").Value = -14196.6671
$ws.Range("H134").Value = 11112550
$ws.Range("I134").Value = 1442.129
$ws.Range("J134").Value = 35715720
$ws.Range("K134").Value = 4326.387
$ws.Range("L134").Value = 107147160
$ws.Range("M134").Value = -1791.387
$ws.Range("N134").Value = -107152230

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3115571
$ws.Range("I4").Value = 9999999
$ws.Range("J4").Value = 1968166.4
$ws.Range("K4").Value = 29999997
$ws.Range("L4").Value = 5904499.199999999
$ws.Range("M4").Value = -29999885
$ws.Range("N4").Value = -5904723.199999999
$ws.Range("H97").Value = 776
$ws.Range("J97").Value = 777.6667
$ws.Range("L97").Value = 2333.0001
$ws.Range("N97").Value = -3325.0001
$ws.Range("H107").Value = 5448.7915
$ws.Range("I107").Value = 497.16666
$ws.Range("J107").Value = 7099.3335
$ws.Range("K107").Value = 1491.49998
$ws.Range("L107").Value = 21298.0005
$ws.Range("M107").Value = 428.5000199999999
$ws.Range("N107").Value = -25138.0005
$ws.Range("H113").Value = 707.6
$ws.Range("I113").Value = 611.75
$ws.Range("J113").Value = 817.1429000000001
$ws.Range("K113").Value = 1835.25
$ws.Range("L113").Value = 2451.4287
$ws.Range("M113").Value = 334.75
$ws.Range("N113").Value = -6791.4287
$ws.Range("H114").Value = 475.73685
$ws.Range("I114").Value = 278
$ws.Range("J114").Value = 747.625
$ws.Range("K114").Value = 834
$ws.Range("L114").Value = 2242.875
$ws.Range("M114").Value = 2420
$ws.Range("N114").Value = -8750.875
$ws.Range("H129").Value = 13441856
$ws.Range("J129").Value = 3969473.5
$ws.Range("L129").Value = 11908420.5
$ws.Range("N129").Value = -11918420.5
$ws.Range("H132").Value = 2123.5
$ws.Range("J132").Value = 2602.5
$ws.Range("L132").Value = 23422.5
$ws.Range("N132").Value = -28482.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 20000
$ws.Range("J15").Value = 20000
$ws.Range("L15").Value = 20000
$ws.Range("N15").Value = -20576
$ws.Range("H81").Value = 20000
$ws.Range("J81").Value = 20000
$ws.Range("L81").Value = 20000
$ws.Range("N81").Value = -21996
$ws.Range("H84").Value = 20000
$ws.Range("J84").Value = 20000
$ws.Range("L84").Value = 60000
$ws.Range("N84").Value = -69984
$ws.Range("H122").Value = 873.1429000000001
$ws.Range("I122").Value = 852
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 2556
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -106
$ws.Range("N122").Value = -7900
$ws.Range("H126").Value = 1982.8
$ws.Range("I126").Value = 1736.2
$ws.Range("J126").Value = 2476
$ws.Range("K126").Value = 5208.6
$ws.Range("L126").Value = 7428
$ws.Range("M126").Value = -2738.6
$ws.Range("N126").Value = -12368
$ws.Range("H132").Value = 5451.5835
$ws.Range("I132").Value = 6360.72
$ws.Range("J132").Value = 3385.3635
$ws.Range("K132").Value = 19082.16
$ws.Range("L132").Value = 10156.0905
$ws.Range("M132").Value = -16552.16
$ws.Range("N132").Value = -15216.0905

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = ""
$ws.Range("N3").Value = ""
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = ""
$ws.Range("N15").Value = ""
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").Value = ""
$ws.Range("H61").Value = 1265.5834
$ws.Range("I61").Value = 1265.5834
$ws.Range("K61").Value = 1265.5834
$ws.Range("M61").Value = -1063.5834
$ws.Range("H113").Value = 1265.5834
$ws.Range("I113").Value = 1265.5834
$ws.Range("K113").Value = 1265.5834
$ws.Range("M113").Value = 904.4166
$ws.Range("H132").Value = 2271.125
$ws.Range("I132").Value = 1748.7727
$ws.Range("K132").Value = 5246.3181
$ws.Range("M132").Value = -2716.3181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1808.1
$ws.Range("I81").Value = 1320.2
$ws.Range("J81").Value = 1905.68
$ws.Range("K81").Value = 2640.4
$ws.Range("L81").Value = 3811.36
$ws.Range("M81").Value = -1579.4
$ws.Range("N81").Value = -5933.360000000001
$ws.Range("H84").Value = 1808.1
$ws.Range("I84").Value = 1320.2
$ws.Range("J84").Value = 1905.68
$ws.Range("K84").Value = 13202
$ws.Range("L84").Value = 19056.8
$ws.Range("M84").Value = -7898
$ws.Range("N84").Value = -29664.8
$ws.Range("H100").Value = 302.53845
$ws.Range("I100").Value = 319
$ws.Range("J100").Value = 247.66667
$ws.Range("K100").Value = 638
$ws.Range("L100").Value = 495.33334
$ws.Range("M100").Value = -97
$ws.Range("N100").Value = -1577.33334
$ws.Range("H107").Value = 355.13043
$ws.Range("I107").Value = 273.64285
$ws.Range("J107").Value = 481.8889
$ws.Range("K107").Value = 820.9285500000001
$ws.Range("L107").Value = 1445.6667
$ws.Range("M107").Value = 1099.07145
$ws.Range("N107").Value = -5285.6667
$ws.Range("H126").Value = 58480570
$ws.Range("I126").Value = 74075160
$ws.Range("K126").Value = 222225480
$ws.Range("M126").Value = -222223010
$ws.Range("H136").Value = 1423.42
$ws.Range("I136").Value = 537.36664
$ws.Range("K136").Value = 1612.09992
$ws.Range("M136").Value = 937.9000800000001
